$d = $word.ActiveDocument

# Names to list under "Docente(s) Responsável(eis)"
$names = @(
    "144651 - Antonio Fernando Sartori",
    "3577649 - Carlos Angelo Nunes",
    "471420 - Carlos Antonio Reis Pereira Baptista",
    "519033 - Carlos Yujiro Shigue",
    "3586455 - Cassius Olivio Figueiredo Terra Ruchert",
    "5840897 - Clodoaldo Saron",
    "5840963 - Daniela Camargo Vernilli",
    "6495737 - Durval Rodrigues Junior",
    "1033242 - Fábio Herbst Florenzano",
    "5983729 - Fernando Vernilli Junior",
    "5009972 - Gilberto Carvalho Coelho",
    "984972 - Hugo Ricardo Zschommler Sandim",
    "1176388 - Luiz Tadeu Fernandes Eleno",
    "7459752 - Maria Ismenia Sodero Toledo Faria",
    "5840622 - Miguel Justino Ribeiro Barboza",
    "2166002 - Sandra Giacomin Schneider",
    "1922320 - Sebastiao Ribeiro",
    "5840793 - Sérgio Schneider"
)

# Find the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Docente(s) Responsável(eis)")) {
        $target = $i
        break
    }
}

# Insert a brand-new paragraph right after the heading; it will host the names.
$headingPara = $d.Paragraphs.Item($target)
$headingPara.Range.InsertParagraphAfter()
$firstIndex = $target + 1

# Build the list: one temporary paragraph per name (each ending with a line
# break for all but the last name), all styled as ListBullet.
for ($n = 0; $n -lt $names.Length; $n++) {
    $para = $d.Paragraphs.Item($firstIndex + $n)
    $para.Range.Style = "ListBullet"
    if ($n -lt $names.Length - 1) {
        $para.Range.InsertAfter($names[$n] + [char]11)
        $para.Range.InsertParagraphAfter()
    } else {
        $para.Range.InsertAfter($names[$n])
    }
}

# Merge all the temporary paragraphs into a single paragraph by deleting the
# paragraph marks between them, which preserves each name (and its line
# break) as its own run instead of collapsing everything into one run.
for ($n = 0; $n -lt $names.Length - 1; $n++) {
    $para = $d.Paragraphs.Item($firstIndex)
    $mark = $d.Range($para.Range.End - 1, $para.Range.End)
    $mark.Delete()
}

Write-Output "done"
